# "aggiornamento fino a 8/12" - extend the daily COVID case table on Sheet1
# with rows for the new dates (A), bringing the series up to serial date
# 44538 (2021-12-08). Columns: A=date, B=nuovi pos., C=somma mobile 7gg.,
# D=somma mobile 7gg. per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstNewRow = 386
$lastNewRow  = 464

# Column A carries a date-time number format + bold/bordered/centered
# style throughout the existing table (see A2:A385). Stamp that same
# formatting onto the new A cells (format only, so no value is copied)
# before writing the real values below.
$ws.Range("A385").Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)

# Build the new rows as a single 2-D array and drop them in with one
# Range.Value assignment.
$newRows = New-Object "object[,]" 79,4
$newRows[0,0] = 44460   # A386 (date serial)
$newRows[0,1] = 1
$newRows[0,2] = 13
$newRows[0,3] = 85.77461071522829
$newRows[1,0] = 44461   # A387 (date serial)
$newRows[1,1] = 0
$newRows[1,2] = 13
$newRows[1,3] = 85.77461071522829
$newRows[2,0] = 44462   # A388 (date serial)
$newRows[2,1] = 0
$newRows[2,2] = 13
$newRows[2,3] = 85.77461071522829
$newRows[3,0] = 44463   # A389 (date serial)
$newRows[3,1] = 2
$newRows[3,2] = 10
$newRows[3,3] = 65.98046978094484
$newRows[4,0] = 44464   # A390 (date serial)
$newRows[4,1] = 0
$newRows[4,2] = 6
$newRows[4,3] = 39.5882818685669
$newRows[5,0] = 44465   # A391 (date serial)
$newRows[5,1] = 0
$newRows[5,2] = 4
$newRows[5,3] = 26.39218791237794
$newRows[6,0] = 44466   # A392 (date serial)
$newRows[6,1] = 0
$newRows[6,2] = 3
$newRows[6,3] = 19.79414093428345
$newRows[7,0] = 44467   # A393 (date serial)
$newRows[7,1] = 2
$newRows[7,2] = 4
$newRows[7,3] = 26.39218791237794
$newRows[8,0] = 44468   # A394 (date serial)
$newRows[8,1] = 0
$newRows[8,2] = 4
$newRows[8,3] = 26.39218791237794
$newRows[9,0] = 44469   # A395 (date serial)
$newRows[9,1] = 0
$newRows[9,2] = 4
$newRows[9,3] = 26.39218791237794
$newRows[10,0] = 44470   # A396 (date serial)
$newRows[10,1] = 0
$newRows[10,2] = 2
$newRows[10,3] = 13.19609395618897
$newRows[11,0] = 44471   # A397 (date serial)
$newRows[11,1] = 0
$newRows[11,2] = 2
$newRows[11,3] = 13.19609395618897
$newRows[12,0] = 44472   # A398 (date serial)
$newRows[12,1] = 0
$newRows[12,2] = 2
$newRows[12,3] = 13.19609395618897
$newRows[13,0] = 44473   # A399 (date serial)
$newRows[13,1] = 0
$newRows[13,2] = 2
$newRows[13,3] = 13.19609395618897
$newRows[14,0] = 44474   # A400 (date serial)
$newRows[14,1] = 0
$newRows[14,2] = 0
$newRows[14,3] = 0
$newRows[15,0] = 44475   # A401 (date serial)
$newRows[15,1] = 0
$newRows[15,2] = 0
$newRows[15,3] = 0
$newRows[16,0] = 44476   # A402 (date serial)
$newRows[16,1] = 1
$newRows[16,2] = 1
$newRows[16,3] = 6.598046978094485
$newRows[17,0] = 44477   # A403 (date serial)
$newRows[17,1] = 0
$newRows[17,2] = 1
$newRows[17,3] = 6.598046978094485
$newRows[18,0] = 44478   # A404 (date serial)
$newRows[18,1] = 0
$newRows[18,2] = 1
$newRows[18,3] = 6.598046978094485
$newRows[19,0] = 44479   # A405 (date serial)
$newRows[19,1] = 0
$newRows[19,2] = 1
$newRows[19,3] = 6.598046978094485
$newRows[20,0] = 44480   # A406 (date serial)
$newRows[20,1] = 0
$newRows[20,2] = 1
$newRows[20,3] = 6.598046978094485
$newRows[21,0] = 44481   # A407 (date serial)
$newRows[21,1] = 0
$newRows[21,2] = 1
$newRows[21,3] = 6.598046978094485
$newRows[22,0] = 44482   # A408 (date serial)
$newRows[22,1] = 0
$newRows[22,2] = 1
$newRows[22,3] = 6.598046978094485
$newRows[23,0] = 44483   # A409 (date serial)
$newRows[23,1] = 0
$newRows[23,2] = 0
$newRows[23,3] = 0
$newRows[24,0] = 44484   # A410 (date serial)
$newRows[24,1] = 0
$newRows[24,2] = 0
$newRows[24,3] = 0
$newRows[25,0] = 44485   # A411 (date serial)
$newRows[25,1] = 0
$newRows[25,2] = 0
$newRows[25,3] = 0
$newRows[26,0] = 44486   # A412 (date serial)
$newRows[26,1] = 0
$newRows[26,2] = 0
$newRows[26,3] = 0
$newRows[27,0] = 44487   # A413 (date serial)
$newRows[27,1] = 0
$newRows[27,2] = 0
$newRows[27,3] = 0
$newRows[28,0] = 44488   # A414 (date serial)
$newRows[28,1] = 0
$newRows[28,2] = 0
$newRows[28,3] = 0
$newRows[29,0] = 44489   # A415 (date serial)
$newRows[29,1] = 0
$newRows[29,2] = 0
$newRows[29,3] = 0
$newRows[30,0] = 44490   # A416 (date serial)
$newRows[30,1] = 0
$newRows[30,2] = 0
$newRows[30,3] = 0
$newRows[31,0] = 44491   # A417 (date serial)
$newRows[31,1] = 0
$newRows[31,2] = 0
$newRows[31,3] = 0
$newRows[32,0] = 44492   # A418 (date serial)
$newRows[32,1] = 0
$newRows[32,2] = 0
$newRows[32,3] = 0
$newRows[33,0] = 44493   # A419 (date serial)
$newRows[33,1] = 0
$newRows[33,2] = 0
$newRows[33,3] = 0
$newRows[34,0] = 44494   # A420 (date serial)
$newRows[34,1] = 0
$newRows[34,2] = 0
$newRows[34,3] = 0
$newRows[35,0] = 44495   # A421 (date serial)
$newRows[35,1] = 0
$newRows[35,2] = 0
$newRows[35,3] = 0
$newRows[36,0] = 44496   # A422 (date serial)
$newRows[36,1] = 0
$newRows[36,2] = 0
$newRows[36,3] = 0
$newRows[37,0] = 44497   # A423 (date serial)
$newRows[37,1] = 0
$newRows[37,2] = 0
$newRows[37,3] = 0
$newRows[38,0] = 44498   # A424 (date serial)
$newRows[38,1] = 0
$newRows[38,2] = 0
$newRows[38,3] = 0
$newRows[39,0] = 44499   # A425 (date serial)
$newRows[39,1] = 0
$newRows[39,2] = 0
$newRows[39,3] = 0
$newRows[40,0] = 44500   # A426 (date serial)
$newRows[40,1] = 1
$newRows[40,2] = 1
$newRows[40,3] = 6.598046978094485
$newRows[41,0] = 44501   # A427 (date serial)
$newRows[41,1] = 0
$newRows[41,2] = 1
$newRows[41,3] = 6.598046978094485
$newRows[42,0] = 44502   # A428 (date serial)
$newRows[42,1] = 0
$newRows[42,2] = 1
$newRows[42,3] = 6.598046978094485
$newRows[43,0] = 44503   # A429 (date serial)
$newRows[43,1] = 0
$newRows[43,2] = 1
$newRows[43,3] = 6.598046978094485
$newRows[44,0] = 44504   # A430 (date serial)
$newRows[44,1] = 0
$newRows[44,2] = 1
$newRows[44,3] = 6.598046978094485
$newRows[45,0] = 44505   # A431 (date serial)
$newRows[45,1] = 3
$newRows[45,2] = 4
$newRows[45,3] = 26.39218791237794
$newRows[46,0] = 44506   # A432 (date serial)
$newRows[46,1] = 0
$newRows[46,2] = 4
$newRows[46,3] = 26.39218791237794
$newRows[47,0] = 44507   # A433 (date serial)
$newRows[47,1] = 3
$newRows[47,2] = 6
$newRows[47,3] = 39.5882818685669
$newRows[48,0] = 44508   # A434 (date serial)
$newRows[48,1] = 0
$newRows[48,2] = 6
$newRows[48,3] = 39.5882818685669
$newRows[49,0] = 44509   # A435 (date serial)
$newRows[49,1] = 2
$newRows[49,2] = 8
$newRows[49,3] = 52.78437582475588
$newRows[50,0] = 44510   # A436 (date serial)
$newRows[50,1] = 0
$newRows[50,2] = 8
$newRows[50,3] = 52.78437582475588
$newRows[51,0] = 44511   # A437 (date serial)
$newRows[51,1] = 2
$newRows[51,2] = 10
$newRows[51,3] = 65.98046978094484
$newRows[52,0] = 44512   # A438 (date serial)
$newRows[52,1] = 0
$newRows[52,2] = 7
$newRows[52,3] = 46.18632884666139
$newRows[53,0] = 44513   # A439 (date serial)
$newRows[53,1] = 0
$newRows[53,2] = 7
$newRows[53,3] = 46.18632884666139
$newRows[54,0] = 44514   # A440 (date serial)
$newRows[54,1] = 2
$newRows[54,2] = 6
$newRows[54,3] = 39.5882818685669
$newRows[55,0] = 44515   # A441 (date serial)
$newRows[55,1] = 0
$newRows[55,2] = 6
$newRows[55,3] = 39.5882818685669
$newRows[56,0] = 44516   # A442 (date serial)
$newRows[56,1] = 3
$newRows[56,2] = 7
$newRows[56,3] = 46.18632884666139
$newRows[57,0] = 44517   # A443 (date serial)
$newRows[57,1] = 0
$newRows[57,2] = 7
$newRows[57,3] = 46.18632884666139
$newRows[58,0] = 44518   # A444 (date serial)
$newRows[58,1] = 0
$newRows[58,2] = 5
$newRows[58,3] = 32.99023489047242
$newRows[59,0] = 44519   # A445 (date serial)
$newRows[59,1] = 10
$newRows[59,2] = 15
$newRows[59,3] = 98.97070467141725
$newRows[60,0] = 44520   # A446 (date serial)
$newRows[60,1] = 0
$newRows[60,2] = 15
$newRows[60,3] = 98.97070467141725
$newRows[61,0] = 44521   # A447 (date serial)
$newRows[61,1] = 3
$newRows[61,2] = 16
$newRows[61,3] = 105.5687516495118
$newRows[62,0] = 44522   # A448 (date serial)
$newRows[62,1] = 4
$newRows[62,2] = 20
$newRows[62,3] = 131.9609395618897
$newRows[63,0] = 44523   # A449 (date serial)
$newRows[63,1] = 19
$newRows[63,2] = 36
$newRows[63,3] = 237.5296912114014
$newRows[64,0] = 44524   # A450 (date serial)
$newRows[64,1] = 16
$newRows[64,2] = 52
$newRows[64,3] = 343.0984428609132
$newRows[65,0] = 44525   # A451 (date serial)
$newRows[65,1] = 1
$newRows[65,2] = 53
$newRows[65,3] = 349.6964898390077
$newRows[66,0] = 44526   # A452 (date serial)
$newRows[66,1] = 41
$newRows[66,2] = 84
$newRows[66,3] = 554.2359461599367
$newRows[67,0] = 44527   # A453 (date serial)
$newRows[67,1] = 5
$newRows[67,2] = 89
$newRows[67,3] = 587.2261810504091
$newRows[68,0] = 44528   # A454 (date serial)
$newRows[68,1] = 24
$newRows[68,2] = 110
$newRows[68,3] = 725.7851675903933
$newRows[69,0] = 44529   # A455 (date serial)
$newRows[69,1] = 0
$newRows[69,2] = 106
$newRows[69,3] = 699.3929796780153
$newRows[70,0] = 44530   # A456 (date serial)
$newRows[70,1] = 25
$newRows[70,2] = 112
$newRows[70,3] = 738.9812615465822
$newRows[71,0] = 44531   # A457 (date serial)
$newRows[71,1] = 1
$newRows[71,2] = 97
$newRows[71,3] = 640.010556875165
$newRows[72,0] = 44532   # A458 (date serial)
$newRows[72,1] = 3
$newRows[72,2] = 99
$newRows[72,3] = 653.206650831354
$newRows[73,0] = 44533   # A459 (date serial)
$newRows[73,1] = 13
$newRows[73,2] = 71
$newRows[73,3] = 468.4613354447084
$newRows[74,0] = 44534   # A460 (date serial)
$newRows[74,1] = 2
$newRows[74,2] = 68
$newRows[74,3] = 448.667194510425
$newRows[75,0] = 44535   # A461 (date serial)
$newRows[75,1] = 22
$newRows[75,2] = 66
$newRows[75,3] = 435.471100554236
$newRows[76,0] = 44536   # A462 (date serial)
$newRows[76,1] = 8
$newRows[76,2] = 74
$newRows[76,3] = 488.2554763789918
$newRows[77,0] = 44537   # A463 (date serial)
$newRows[77,1] = 8
$newRows[77,2] = 57
$newRows[77,3] = 376.0886777513856
$newRows[78,0] = 44538   # A464 (date serial)
$newRows[78,1] = 1
$newRows[78,2] = 57
$newRows[78,3] = 376.0886777513856

$ws.Range("A" + $firstNewRow + ":D" + $lastNewRow).Value = $newRows

Write-Output ("Sheet1 extended to row 464: A" + $firstNewRow + ":D" + $lastNewRow + " filled (through date serial 44538).")
